$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text (inlineStr) in the
# original workbook, so each new value is force-typed as text (leading
# apostrophe) and then restyled back to the sheet's default so no stray
# number-format/style is introduced.
$ws.Range("D2").Value = "'274.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.208"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06177"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.575"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.514"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'6.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8231"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1640"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.08212"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03429"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03144"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09132"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.764"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Value = "'0.04695"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.006454"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006138"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'3.724"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'2.317"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'0.01388"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3277"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Value = "'0.0002737"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.04745"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.005500"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'40CEJICEJIBestin24h"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007024"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'41KickTokenKICK"
$ws.Range("E42").Style = "Normal"
$ws.Range("D44").Value = "'0.01024"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006564"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.7233"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Value = "'0.00001900"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.01240"
$ws.Range("D50").Style = "Normal"
